$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 901, shifting existing rows 901-936 down to 903-938
$ws.Rows.Item(901).Insert()
$ws.Rows.Item(901).Insert()

# Populate new row 901 (new price record)
$ws.Cells.Item(901, 1).Value = 9
$ws.Cells.Item(901, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(901, 3).Value = "Metropolitana"
$ws.Cells.Item(901, 4).Value = 45267
$ws.Cells.Item(901, 5).Value = 13
$ws.Cells.Item(901, 6).Value = "Fruta"
$ws.Cells.Item(901, 7).Value = 100109
$ws.Cells.Item(901, 8).Value = "Uva"
$ws.Cells.Item(901, 9).Value = 100109001
$ws.Cells.Item(901, 10).Value = "Uva"
$ws.Cells.Item(901, 11).Value = "Flame Seedless"
$ws.Cells.Item(901, 12).Value = "Primera"
$ws.Cells.Item(901, 13).Value = 240
$ws.Cells.Item(901, 14).Value = 20000
$ws.Cells.Item(901, 15).Value = 20000
$ws.Cells.Item(901, 16).Value = 20000
$ws.Cells.Item(901, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(901, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(901, 19).Value = 2000
$ws.Cells.Item(901, 20).Value = 10

# Populate new row 902 (new price record)
$ws.Cells.Item(902, 1).Value = 9
$ws.Cells.Item(902, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(902, 3).Value = "Metropolitana"
$ws.Cells.Item(902, 4).Value = 45267
$ws.Cells.Item(902, 5).Value = 13
$ws.Cells.Item(902, 6).Value = "Fruta"
$ws.Cells.Item(902, 7).Value = 100109
$ws.Cells.Item(902, 8).Value = "Uva"
$ws.Cells.Item(902, 9).Value = 100109001
$ws.Cells.Item(902, 10).Value = "Uva"
$ws.Cells.Item(902, 11).Value = "Superior Seedless"
$ws.Cells.Item(902, 12).Value = "Primera"
$ws.Cells.Item(902, 13).Value = 150
$ws.Cells.Item(902, 14).Value = 25000
$ws.Cells.Item(902, 15).Value = 25000
$ws.Cells.Item(902, 16).Value = 25000
$ws.Cells.Item(902, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(902, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(902, 19).Value = 2500
$ws.Cells.Item(902, 20).Value = 10
